$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump "Last Updated" timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "30 Oct 2025, 10:05 AM"

# --- "distance from Dma50" sheet: refresh values, swap two rows ---
$ws = $wb.Worksheets.Item("distance from Dma50")

$ws.Range("C2").Value  = 9.3697
$ws.Range("C3").Value  = 7.3204
$ws.Range("C4").Value  = 5.7935
$ws.Range("C5").Value  = 5.0421
$ws.Range("C6").Value  = 4.9843
$ws.Range("C7").Value  = 4.8829
$ws.Range("C8").Value  = 4.298
$ws.Range("C9").Value  = 4.247
$ws.Range("C10").Value = 3.4564
$ws.Range("C11").Value = 3.4146
$ws.Range("C12").Value = 3.2672
$ws.Range("C13").Value = 3.0616
$ws.Range("C14").Value = 3.0453
$ws.Range("C15").Value = 2.9329
$ws.Range("C16").Value = 2.8979
$ws.Range("C17").Value = 2.713
$ws.Range("C18").Value = 2.5331
$ws.Range("C19").Value = 2.2569
$ws.Range("C20").Value = 2.1293
$ws.Range("C21").Value = 2.074

# Rows 22/23 swap stock names (CNXIT <-> NIFTYCONSUMPTION) along with new values
$ws.Range("B22").Value = "CNXIT"
$ws.Range("C22").Value = 1.2806
$ws.Range("B23").Value = "NIFTYCONSUMPTION"
$ws.Range("C23").Value = 1.2682

$ws.Range("C24").Value = 0.9797
$ws.Range("C25").Value = 0.832
$ws.Range("C26").Value = 0.7733
$ws.Range("C27").Value = 0.6065
$ws.Range("C28").Value = 0.1789
$ws.Range("C29").Value = -0.2561
$ws.Range("C30").Value = -2.0251
